$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.647.96"
$ws.Range("E2").Value = "  +1.46%  "

# Row 3
$ws.Range("D3").Value = "1.802.45"
$ws.Range("E3").Value = "  +1.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'227.45"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6
$ws.Range("E6").Value = "  +1.93%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("D8").Value = "'32.88"
$ws.Range("E8").Value = "  +4.19%  "

# Row 9
$ws.Range("E9").Value = "  +1.96%  "

# Row 10
$ws.Range("D10").Value = "'0.0696"

# Row 12
$ws.Range("D12").Value = "2.064.56"
$ws.Range("E12").Value = "  +1.45%  "

# Row 13
$ws.Range("D13").Value = "'11.21"
$ws.Range("E13").Value = "  +2.77%  "

# Row 14
$ws.Range("D14").Value = "1.799.79"
$ws.Range("E14").Value = "  +1.31%  "

# Row 15
$ws.Range("D15").Value = "'0.640"
$ws.Range("E15").Value = "  +3.06%  "

# Row 16
$ws.Range("D16").Value = "34.640.02"
$ws.Range("E16").Value = "  +1.56%  "

# Row 17
$ws.Range("E17").Value = "  +3.75%  "

# Row 18
$ws.Range("D18").Value = "'69.03"
$ws.Range("E18").Value = "  +1.81%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0804"
$ws.Range("E19").Value = "  +0.71%  "

# Row 20
$ws.Range("D20").Value = "'247.46"
$ws.Range("E20").Value = "  +0.85%  "

# Row 21
$ws.Range("D21").Value = "'11.35"
$ws.Range("E21").Value = "  +3.51%  "

# Row 22
$ws.Range("E22").Value = "  -0.35%  "

# Row 23
$ws.Range("D23").Value = "'4.18"
$ws.Range("E23").Value = "  +2.38%  "

# Row 24
$ws.Range("D24").Value = "'171.63"
$ws.Range("E24").Value = "  +5.62%  "

# Row 25
$ws.Range("E25").Value = "  +2.03%  "

# Row 26
$ws.Range("D26").Value = "'7.34"
$ws.Range("E26").Value = "  +2.18%  "

# Row 27
$ws.Range("D27").Value = "'16.63"
$ws.Range("E27").Value = "  +2.20%  "

# Row 28
$ws.Range("D28").Value = "'0.116"
$ws.Range("E28").Value = "  +2.27%  "

# Row 29
$ws.Range("E29").Value = "  -0.31%  "

# Row 30
$ws.Range("D30").Value = "'4.11"
$ws.Range("E30").Value = "  +11.17%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0527"
$ws.Range("E31").Value = "  +1.44%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  +1.10%  "

# Row 33
$ws.Range("E33").Value = "  +2.45%  "

# Row 34
$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  +3.02%  "

# Row 35
$ws.Range("D35").Value = "1.435.13"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("D36").Value = "'2.58"
$ws.Range("E36").Value = "  +7.32%  "

# Row 37
$ws.Range("D37").Value = "'0.678"
$ws.Range("E37").Value = "  +3.13%  "

# Row 38
$ws.Range("D38").Value = "'1.07"
$ws.Range("E38").Value = "  +2.74%  "

# Row 39
$ws.Range("D39").Value = "'0.0191"
$ws.Range("E39").Value = "  +0.67%  "

# Row 40
$ws.Range("D40").Value = "'85.10"
$ws.Range("E40").Value = "  +6.26%  "

# Row 41
$ws.Range("E41").Value = "  +3.70%  "

# Row 42
$ws.Range("E42").Value = "  +1.45%  "

# Row 43
$ws.Range("E43").Value = "  +3.58%  "

# Row 44
$ws.Range("D44").Value = "'13.87"
$ws.Range("E44").Value = "  +2.97%  "

# Row 45
$ws.Range("D45").Value = "'0.0526"
$ws.Range("E45").Value = "  +2.97%  "

# Row 46
$ws.Range("D46").Value = "'6.12"
$ws.Range("E46").Value = "  +0.72%  "

# Row 47
$ws.Range("E47").Value = "  +0.49%  "

# Row 48
$ws.Range("D48").Value = "1.965.32"
$ws.Range("E48").Value = "  +1.40%  "

# Row 49
$ws.Range("D49").Value = "'105.67"
$ws.Range("E49").Value = "  +1.42%  "

# Row 50
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.31%  "

# Row 51
$ws.Range("E51").Value = "  -5.05%  "
